# feat: add 2022-Q1 data
#
# Source workbook has 2 sheets: "2021-Q4" (fund-holding detail) and "总计"
# (summary). This adds a new "2022-Q1" fund-holding detail sheet (between
# "2021-Q4" and "总计") and prepends a 2022-Q1 row to the "总计" summary
# sheet (pushing the existing 2021-Q4 row down one row).

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetOrig = $wb.Worksheets.Item("总计")

# --- 1. Insert the new "2022-Q1" detail sheet, positioned right before
#        "总计" (so sheet order becomes 2021-Q4, 2022-Q1, 总计). Duplicate
#        the 2021-Q4 sheet so layout/styles (header row + index column
#        bold+border style, column widths, etc.) come along for free, then
#        overwrite the cell values below.
$q4Sheet.Copy($totalSheetOrig)
$q1Sheet = $wb.Worksheets.Item(2)
$q1Sheet.Name = "2022-Q1"

# Worksheet handles returned earlier are positional, and the sheet that
# used to sit at "总计"'s slot has just shifted one place to the right, so
# re-resolve "总计" by name now that the insert above is done.
$totalSheet = $wb.Worksheets.Item("总计")

# Fund rows for 2022-Q1 (code, name, scale, stock position, position pct,
# market value, position rank). Every column except A (index) and H (rank)
# is stored as text in the source data, even though it looks numeric, so
# each value is assigned with a leading "'" to force text and then the
# style is reset to Normal so no stray NumberFormat/quote-prefix style is
# left behind on the cell (matches the unstyled inlineStr cells in the
# target file).
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $q1Sheet.Range("B2") "519702"
Set-TextCell $q1Sheet.Range("C2") "交银趋势优先混合"
Set-TextCell $q1Sheet.Range("D2") "108.29"
Set-TextCell $q1Sheet.Range("E2") "71.40"
Set-TextCell $q1Sheet.Range("F2") "2.35"
Set-TextCell $q1Sheet.Range("G2") "2.5448"
$q1Sheet.Range("H2").Value = 5

Set-TextCell $q1Sheet.Range("B3") "006143"
Set-TextCell $q1Sheet.Range("C3") "恒生前海中证质量成长低波动指数A"
Set-TextCell $q1Sheet.Range("D3") "0.06"
Set-TextCell $q1Sheet.Range("E3") "94.34"
Set-TextCell $q1Sheet.Range("F3") "2.29"
Set-TextCell $q1Sheet.Range("G3") "0.0014"
$q1Sheet.Range("H3").Value = 8

Set-TextCell $q1Sheet.Range("B4") "006144"
Set-TextCell $q1Sheet.Range("C4") "恒生前海中证质量成长低波动指数C"
Set-TextCell $q1Sheet.Range("D4") "0.01"
Set-TextCell $q1Sheet.Range("E4") "94.34"
Set-TextCell $q1Sheet.Range("F4") "2.29"
Set-TextCell $q1Sheet.Range("G4") "0.0002"
$q1Sheet.Range("H4").Value = 8

# --- 2. Update the "总计" summary sheet: push the existing 2021-Q4 row
#        (row 2) down to row 3, then write the new 2022-Q1 totals into row
#        2. Copy-to-destination (rather than cut/insert) preserves the
#        bold/border style on column A exactly.
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 2.55

$totalSheet.Range("A3").Value = 1
